$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.911.04'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.622.75'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '213.63'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('E10').Value = '  -6.46%  '
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').Value = '1.848.79'
$ws.Range('D13').Value = '1.619.28'
$ws.Range('E13').Value = '  -2.28%  '
$ws.Range('D14').Value = '4.17'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '25.906.96'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '61.12'
$ws.Range('E17').Value = '  -3.55%  '
$ws.Range('D18').Value = '0.0₃0733'
$ws.Range('E18').Value = '  -3.90%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '191.91'
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '0.133'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').Value = '143.86'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -2.91%  '
$ws.Range('D28').Value = '6.71'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('E29').Value = '  -2.46%  '
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('D32').Value = '3.11'
$ws.Range('E32').Value = '  -4.35%  '
$ws.Range('E33').Value = '  -5.55%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').Value = '1.114.94'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('E37').Value = '  -6.52%  '
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').Value = '97.94'
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('D42').Value = '0.765'
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('D43').Value = '1.758.86'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  -5.83%  '
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('D47').Value = '54.26'
$ws.Range('E47').Value = '  -3.99%  '
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '7.44'
$ws.Range('E51').Value = '  -3.63%  '
